# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 388
$ws1.Range("F5").Value = 432
$ws1.Range("F7").Value = 2441
$ws1.Range("F8").Value = 422
$ws1.Range("F9").Value = 6401
$ws1.Range("F12").Value = 25

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 388
$ws4.Range("F5").Value = 432
$ws4.Range("F9").Value = 2441
$ws4.Range("F10").Value = 422
$ws4.Range("F11").Value = 6401
$ws4.Range("F15").Value = 25
